{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nfor (let i = 0; i < count; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n  p.inlinePictures.load(\"items\");\n  p.parentTableOrNullObject.load(\"isNullObject\");\n}\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < count; i++) {\n  const p = paragraphs.items[i];\n  const inTable = !p.parentTableOrNullObject.isNullObject;\n  const isEmptyText = p.text.trim().length === 0;\n  const hasPicture = p.inlinePictures.items.length > 0;\n  const isLastBodyParagraph = i === count - 1;\n  // Remove the \"screenshot\" paragraphs that hold nothing but an inline\n  // picture, and the empty spacer paragraphs (no text, no picture, not\n  // the very last paragraph of the document, not inside a table) that\n  // Word leaves right after a table.\n  if (!inTable && (hasPicture || isEmptyText) && !isLastBodyParagraph) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$wdWithInTable = 12\n$docEnd = $d.Content.End\n\n# Collect the ranges that must be removed:\n#   1) every paragraph that is made up solely of an inline picture\n#      (the three \"screenshot\" paragraphs under the headings), and\n#   2) every empty spacer paragraph (pPr with only spacing before=40,\n#      i.e. 2pt) that Word leaves right after a table.\n# We gather Range.Start/Range.End pairs first and then delete them\n# back-to-front so earlier deletions don't shift the positions of the\n# ones still pending.\n\n$ranges = New-Object System.Collections.ArrayList\n\nfor ($i = 1; $i -le $d.InlineShapes.Count; $i++) {\n    $shape = $d.InlineShapes.Item($i)\n    $para = $shape.Range.Paragraphs.Item(1)\n    [void]$ranges.Add(@{Start = $para.Range.Start; End = $para.Range.End})\n}\n\n$n = $d.Paragraphs.Count\nfor ($i = 1; $i -le $n; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $inTable = $p.Range.Information($wdWithInTable)\n    $isEmpty = ($p.Range.InlineShapes.Count -eq 0) -and ($p.Range.Text.Length -eq 1)\n    $isLastParagraph = ($p.Range.End -eq $docEnd)\n    if (-not $inTable -and $isEmpty -and ($p.SpaceBefore -eq 2) -and -not $isLastParagraph) {\n        [void]$ranges.Add(@{Start = $p.Range.Start; End = $p.Range.End})\n    }\n}\n\n$sorted = $ranges | Sort-Object -Property Start -Descending\n\nforeach ($r in $sorted) {\n    $rng = $d.Range($r.Start, $r.End)\n    $rng.Delete()\n}\n"}
